$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Timing_Category")

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "Hindcast"
$ws.Range("C11").Value = "Use the model to estimate unobserved past conditions"

$wb.Names.Item("Timing_Category").RefersTo = '=Timing_Category!$A$1:$C$11'
